$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "86.806.92"
$ws.Range("E2").Value = "  -1.89%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.149.10"
$ws.Range("E3").Value = "  -5.33%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.13%  "

# Row 5 - Solana
$ws.Range("D5").Value = "'205.48"
$ws.Range("E5").Value = "  -5.46%  "

# Row 6 - BNB
$ws.Range("D6").Value = "'604.06"
$ws.Range("E6").Value = "  -6.95%  "

# Row 7 - Dogecoin
$ws.Range("D7").Value = "'0.361"
$ws.Range("E7").Value = "  -9.10%  "

# Row 8 - XRP
$ws.Range("D8").Value = "'0.649"
$ws.Range("E8").Value = "  +8.10%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  +0.07%  "

# Row 10 - LidoStakedEther
$ws.Range("D10").Value = "3.148.72"
$ws.Range("E10").Value = "  -5.18%  "

# Row 11 - Cardano
$ws.Range("D11").Value = "'0.530"
$ws.Range("E11").Value = "  -9.36%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +5.35%  "

# Row 13 - ShibaInu
$ws.Range("D13").Value = "'0.0000240"
$ws.Range("E13").Value = "  -17.54%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.737.00"
$ws.Range("E14").Value = "  -5.07%  "

# Row 15 - Toncoin
$ws.Range("D15").Value = "'5.22"
$ws.Range("E15").Value = "  -5.24%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "86.567.39"
$ws.Range("E16").Value = "  -2.00%  "

# Row 17 - Avalanche
$ws.Range("D17").Value = "'31.75"
$ws.Range("E17").Value = "  -10.41%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.171.82"
$ws.Range("E18").Value = "  -4.68%  "

# Row 19 - SuiNetwork
$ws.Range("D19").Value = "'2.93"
$ws.Range("E19").Value = "  -5.60%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "'13.24"
$ws.Range("E20").Value = "  -9.30%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "'408.72"
$ws.Range("E21").Value = "  -10.05%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'8.38"
$ws.Range("E22").Value = "  -12.70%  "

# Row 23 - Polkadot
$ws.Range("D23").Value = "'4.99"
$ws.Range("E23").Value = "  -7.89%  "

# Row 24 - NEARProtocol
$ws.Range("D24").Value = "'5.08"
$ws.Range("E24").Value = "  -8.38%  "

# Row 25 - Aptos
$ws.Range("D25").Value = "'11.47"
$ws.Range("E25").Value = "  -11.04%  "

# Row 26 - WrappedeETH
$ws.Range("D26").Value = "3.344.72"
$ws.Range("E26").Value = "  -4.76%  "

# Row 27 - Litecoin
$ws.Range("D27").Value = "'72.88"
$ws.Range("E27").Value = "  -6.79%  "

# Row 28 - PEPE
$ws.Range("D28").Value = "'0.0000128"
$ws.Range("E28").Value = "  -4.03%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  -0.01%  "

# Row 30 - Cronos
$ws.Range("D30").Value = "'0.159"
$ws.Range("E30").Value = "  -25.01%  "

# Row 31 - Binance-PegBSC-USD
$ws.Range("D31").Value = "'0.996"
$ws.Range("E31").Value = "  -0.46%  "

# Row 32 - Bittensor
$ws.Range("D32").Value = "'533.67"
$ws.Range("E32").Value = "  -10.36%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Value = "'8.20"
$ws.Range("E33").Value = "  -12.33%  "

# Row 34 - PancakeSwap
$ws.Range("D34").Value = "'1.84"
$ws.Range("E34").Value = "  -13.11%  "

# Row 35 - Fetch.AI
$ws.Range("D35").Value = "'1.26"
$ws.Range("E35").Value = "  -20.88%  "

# Row 36 - RenderToken
$ws.Range("D36").Value = "'6.51"
$ws.Range("E36").Value = "  -9.28%  "

# Row 37 - Kaspa
$ws.Range("D37").Value = "'0.132"
$ws.Range("E37").Value = "  -6.19%  "

# Row 38 - EthereumClassic
$ws.Range("D38").Value = "'21.50"
$ws.Range("E38").Value = "  -6.90%  "

# Row 39 - WhiteBITCoin
$ws.Range("D39").Value = "'21.78"
$ws.Range("E39").Value = "  -0.24%  "

# Row 40 - FirstDigitalUSD
$ws.Range("D40").Value = "'0.997"
$ws.Range("E40").Value = "  +0.05%  "

# Row 41 - dogwifhat
$ws.Range("D41").Value = "'2.96"
$ws.Range("E41").Value = "  -5.36%  "

# Row 42 - USDe
$ws.Range("E42").Value = "  -0.02%  "

# Row 43 - PolygonEcosystemToken
$ws.Range("D43").Value = "'0.370"
$ws.Range("E43").Value = "  -11.65%  "

# Row 44 - Stacks
$ws.Range("D44").Value = "'1.88"
$ws.Range("E44").Value = "  -12.80%  "

# Row 45 - Monero
$ws.Range("D45").Value = "'149.82"
$ws.Range("E45").Value = "  -5.24%  "

# Rows 46/47 swap places: OKB moves to 46, Aave moves to 47 (A column rank stays put)
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").Value = "'43.09"
$ws.Range("E46").Value = "  -6.22%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'170.51"
$ws.Range("E47").Value = "  -9.06%  "

# Row 48 - Stellar
$ws.Range("D48").Value = "'0.126"
$ws.Range("E48").Value = "  +8.67%  "

# Row 49 - ImmutableX
$ws.Range("D49").Value = "'1.24"
$ws.Range("E49").Value = "  -13.64%  "

# Row 50 - Filecoin
$ws.Range("D50").Value = "'3.91"
$ws.Range("E50").Value = "  -11.73%  "

# Row 51 - ARBITRUM
$ws.Range("D51").Value = "'0.581"
$ws.Range("E51").Value = "  -11.69%  "
